$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, shifting existing rows down
$ws.Rows.Item(6).Insert()

# Copy formatting from row 8 (which now has the "banded" style matching what row 6 should have)
$ws.Range("B8:G8").Copy()
$ws.Range("B6:G6").PasteSpecial(-4122)

# Populate the new row's values (August 2025 data)
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Ago."
$ws.Range("D6").Value = 360.685
$ws.Range("E6").Value = 31756.39
$ws.Range("F6").Value = 4424.647
$ws.Range("G6").Value = 149.563

# Resize table / autofilter to include new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:G97"))

# Update the "Actualizacion" text cell
$ws.Range("B98").Value = "Actualización: Agosto 2025."
